$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.717.77'
$ws.Cells.Item(2, 5).Value = '  -2.70%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.095.00'
$ws.Cells.Item(3, 5).Value = '  -1.98%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.011'
$ws.Cells.Item(4, 5).Value = '  +0.17%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '344.10'
$ws.Cells.Item(5, 5).Value = '  -2.29%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '1.009'
$ws.Cells.Item(6, 5).Value = '  +0.21%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.5160'
$ws.Cells.Item(7, 5).Value = '  -1.88%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.4377'
$ws.Cells.Item(8, 5).Value = '  -3.95%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '53.05'
$ws.Cells.Item(9, 5).Value = '  -0.94%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +0.88%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -1.92%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '24.84'
$ws.Cells.Item(12, 5).Value = '  -2.19%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '2.098.86'
$ws.Cells.Item(13, 5).Value = '  -1.91%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '8.324'
$ws.Cells.Item(14, 5).Value = '  +1.57%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '6.740'
$ws.Cells.Item(15, 5).Value = '  -2.30%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '99.29'
$ws.Cells.Item(16, 5).Value = '  -2.95%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.00001150'
$ws.Cells.Item(17, 5).Value = '  -2.08%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '1.010'
$ws.Cells.Item(18, 5).Value = '  +0.16%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '20.75'
$ws.Cells.Item(19, 5).Value = '  +1.20%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.06664'
$ws.Cells.Item(20, 5).Value = '  -0.84%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '1.010'
$ws.Cells.Item(21, 5).Value = '  +0.26%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.182'
$ws.Cells.Item(22, 5).Value = '  -3.08%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '29.755.77'
$ws.Cells.Item(23, 5).Value = '  -2.87%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '12.49'
$ws.Cells.Item(24, 5).Value = '  -3.19%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.322'
$ws.Cells.Item(25, 5).Value = '  -2.92%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '2.347.17'
$ws.Cells.Item(26, 5).Value = '  -1.83%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '21.91'
$ws.Cells.Item(27, 5).Value = '  -2.70%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.515'
$ws.Cells.Item(28, 5).Value = '  -4.96%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Monero'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '161.37'
$ws.Cells.Item(29, 5).Value = '  -2.18%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '133.02'
$ws.Cells.Item(30, 5).Value = '  -2.23%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.129'
$ws.Cells.Item(31, 5).Value = '  -7.77%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -2.94%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.648'
$ws.Cells.Item(33, 5).Value = '  -2.48%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '6.149'
$ws.Cells.Item(34, 5).Value = '  -4.18%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '3.934'
$ws.Cells.Item(35, 5).Value = '  -2.56%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '6.226'
$ws.Cells.Item(36, 5).Value = '  +1.42%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -3.34%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.02570'
$ws.Cells.Item(38, 5).Value = '  -3.07%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.06685'
$ws.Cells.Item(39, 5).Value = '  -4.20%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '1.338'
$ws.Cells.Item(40, 5).Value = '  +4.50%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -2.63%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.6862'
$ws.Cells.Item(42, 5).Value = '  -2.00%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.2224'
$ws.Cells.Item(43, 5).Value = '  -4.88%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.6678'
$ws.Cells.Item(44, 5).Value = '  +2.55%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '14.26'
$ws.Cells.Item(45, 5).Value = '  -3.74%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '2.312'
$ws.Cells.Item(46, 5).Value = '  -2.08%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.00000000360'
$ws.Cells.Item(47, 5).Value = '  -2.80%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.624'
$ws.Cells.Item(48, 5).Value = '  -3.44%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.216'
$ws.Cells.Item(49, 5).Value = '  -2.78%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -2.75%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'WOONetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.3260'
$ws.Cells.Item(51, 5).Value = '  -1.67%  '
